$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 387.83334
$ws.Range("I33").Value = 186.625
$ws.Range("J33").Value = 790.25
$ws.Range("K33").Value = 186.625
$ws.Range("L33").Value = 790.25
$ws.Range("M33").Value = 42.375
$ws.Range("N33").Value = -1248.25

$ws.Range("H64").Value = 3174.5454
$ws.Range("I64").Value = 3100
$ws.Range("J64").Value = 3373.3333
$ws.Range("K64").Value = 3100
$ws.Range("L64").Value = 3373.3333
$ws.Range("M64").Value = -2852
$ws.Range("N64").Value = -3869.3333

$ws.Range("H67").Value = 3174.5454
$ws.Range("I67").Value = 3100
$ws.Range("J67").Value = 3373.3333
$ws.Range("K67").Value = 3100
$ws.Range("L67").Value = 3373.3333
$ws.Range("M67").Value = -2242
$ws.Range("N67").Value = -5089.3333

$ws.Range("H137").Value = 1124.7715
$ws.Range("I137").Value = 937.13794
$ws.Range("J137").Value = 2031.6666
$ws.Range("K137").Value = 2811.41382
$ws.Range("L137").Value = 6094.9998
$ws.Range("M137").Value = -261.4138199999998
$ws.Range("N137").Value = -11194.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4166.222
$ws.Range("I45").Value = 3668
$ws.Range("J45").Value = 4664.4443
$ws.Range("K45").Value = 3668
$ws.Range("L45").Value = 4664.4443
$ws.Range("M45").Value = -3291
$ws.Range("N45").Value = -5418.4443

$ws.Range("H61").Value = 4058.8
$ws.Range("I61").Value = 3733.4546
$ws.Range("J61").Value = 4953.5
$ws.Range("K61").Value = 3733.4546
$ws.Range("L61").Value = 4953.5
$ws.Range("M61").Value = -3521.4546
$ws.Range("N61").Value = -5377.5

$ws.Range("H102").Value = 50001536
$ws.Range("I102").Value = 1589.7858
$ws.Range("J102").Value = 166668080
$ws.Range("K102").Value = 1589.7858
$ws.Range("L102").Value = 166668080
$ws.Range("M102").Value = 32.21419999999989
$ws.Range("N102").Value = -166671324

$ws.Range("H108").Value = 29266.666
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 29266.666
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 29266.666
$ws.Range("N108").Value = -36946.666

$ws.Range("H132").Value = 2594
$ws.Range("I132").Value = 1843.875
$ws.Range("J132").Value = 3344.125
$ws.Range("K132").Value = 5531.625
$ws.Range("L132").Value = 10032.375
$ws.Range("M132").Value = -3001.625
$ws.Range("N132").Value = -15092.375

$ws.Range("H136").Value = 4058.8
$ws.Range("I136").Value = 3733.4546
$ws.Range("J136").Value = 4953.5
$ws.Range("K136").Value = 11200.3638
$ws.Range("L136").Value = 14860.5
$ws.Range("M136").Value = -8650.363799999999
$ws.Range("N136").Value = -19960.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2360.524
$ws.Range("I99").Value = 1420.8334
$ws.Range("J99").Value = 3613.4443
$ws.Range("K99").Value = 1420.8334
$ws.Range("L99").Value = 3613.4443
$ws.Range("M99").Value = 77.16660000000002
$ws.Range("N99").Value = -6609.4443

$ws.Range("H105").Value = 3535.0715
$ws.Range("I105").Value = 2182.8572
$ws.Range("J105").Value = 4887.2856
$ws.Range("K105").Value = 2182.8572
$ws.Range("L105").Value = 4887.2856
$ws.Range("M105").Value = -435.8571999999999
$ws.Range("N105").Value = -8381.285599999999

$ws.Range("H107").Value = 1701.68
$ws.Range("I107").Value = 1618
$ws.Range("J107").Value = 1966.6666
$ws.Range("K107").Value = 1618
$ws.Range("L107").Value = 1966.6666
$ws.Range("M107").Value = 302
$ws.Range("N107").Value = -5806.6666

$ws.Range("H118").Value = 16890.334
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 16890.334
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 16890.334
$ws.Range("N118").Value = -20204.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 884.2308
$ws.Range("I22").Value = 885.9091
$ws.Range("J22").Value = 875
$ws.Range("K22").Value = 885.9091
$ws.Range("L22").Value = 875
$ws.Range("M22").Value = -535.9091
$ws.Range("N22").Value = -1575

$ws.Range("H31").Value = 16395367
$ws.Range("I31").Value = 40001320
$ws.Range("J31").Value = 2344.25
$ws.Range("K31").Value = 40001320
$ws.Range("L31").Value = 2344.25
$ws.Range("M31").Value = -40001025
$ws.Range("N31").Value = -2934.25

$ws.Range("H34").Value = 16395367
$ws.Range("I34").Value = 40001320
$ws.Range("J34").Value = 2344.25
$ws.Range("K34").Value = 40001320
$ws.Range("L34").Value = 2344.25
$ws.Range("M34").Value = -40001118
$ws.Range("N34").Value = -2748.25

$ws.Range("H62").Value = 2745.1177
$ws.Range("I62").Value = 2735.4375
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 2735.4375
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -2111.4375
$ws.Range("N62").Value = -4148

$ws.Range("H65").Value = 2745.1177
$ws.Range("I65").Value = 2735.4375
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 13677.1875
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -10557.1875
$ws.Range("N65").Value = -20740

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 27994.666
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 27994.666
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 27994.666
$ws.Range("N15").Value = -28570.666

$ws.Range("H45").Value = 13227.083
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 13227.083
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 13227.083
$ws.Range("N45").Value = -14345.083

$ws.Range("H81").Value = 27994.666
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 27994.666
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 27994.666
$ws.Range("N81").Value = -29990.666

$ws.Range("H84").Value = 27994.666
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 27994.666
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 83983.99800000001
$ws.Range("N84").Value = -93967.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1511
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 1977.5
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 1977.5
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -2353.5

$ws.Range("H61").Value = 2225.2307
$ws.Range("I61").Value = 1265.7273
$ws.Range("J61").Value = 7502.5
$ws.Range("K61").Value = 1265.7273
$ws.Range("L61").Value = 7502.5
$ws.Range("M61").Value = -1063.7273
$ws.Range("N61").Value = -7906.5

$ws.Range("H80").Value = 33500
$ws.Range("I80").Value = 30000
$ws.Range("J80").Value = 37000
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 37000
$ws.Range("M80").Value = -28877
$ws.Range("N80").Value = -39246

$ws.Range("H83").Value = 33500
$ws.Range("I83").Value = 30000
$ws.Range("J83").Value = 37000
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 111000
$ws.Range("M83").Value = -84384
$ws.Range("N83").Value = -122232

$ws.Range("H113").Value = 2225.2307
$ws.Range("I113").Value = 1265.7273
$ws.Range("J113").Value = 7502.5
$ws.Range("K113").Value = 1265.7273
$ws.Range("L113").Value = 7502.5
$ws.Range("M113").Value = 904.2727
$ws.Range("N113").Value = -11842.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3470.7778
$ws.Range("I132").Value = 3697.923
$ws.Range("J132").Value = 2880.2
$ws.Range("K132").Value = 11093.769
$ws.Range("L132").Value = 8640.599999999999
$ws.Range("M132").Value = -8563.769
$ws.Range("N132").Value = -13700.6
